# Update automatico via Actualizar 03-11-2021 12-32-08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timestamp values (Excel date serials), one "generation" shifted down:
#  rows 2-15  -> newest timestamp
#  rows 16-29 -> previous timestamp (was rows 2-15's old value)
#  rows 30-43 -> previous-previous timestamp (was rows 16-29's old value)

$newest = 44266.52208147047
$mid    = 44266.50069510417
$old    = 44266.47930578703

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $mid
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $old
}
